# Apply updates described by the diff to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Update existing odds values in row 2 (Finnish Veikkausliiga)
# ---------------------------------------------------------------------
$ws.Range("G2").Value = 2.14
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 3.65
$ws.Range("L2").Value = 1.22
$ws.Range("N2").Value = 6
$ws.Range("V2").Value = 1.38
$ws.Range("W2").Value = 1.87
$ws.Range("X2").Value = 32
$ws.Range("Y2").Value = 980
$ws.Range("Z2").Value = 38
$ws.Range("AB2").Value = 16.5
$ws.Range("AE2").Value = 980
$ws.Range("AF2").Value = 21
$ws.Range("AG2").Value = 14
$ws.Range("AH2").Value = 18.5
$ws.Range("AJ2").Value = 980

# ---------------------------------------------------------------------
# 2) Update existing odds values in row 3 (Danish 1st Division)
# ---------------------------------------------------------------------
$ws.Range("J3").Value = 3.75
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 2.3
$ws.Range("P3").Value = 2.3
$ws.Range("Q3").Value = 1.62
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 1.03
$ws.Range("U3").Value = 2.24
$ws.Range("X3").Value = 1000
$ws.Range("AA3").Value = 95
$ws.Range("AK3").Value = 29
$ws.Range("AL3").Value = 44
$ws.Range("AM3").Value = 95
$ws.Range("AN3").Value = 16

# ---------------------------------------------------------------------
# 3) Update existing odds values in row 4 (Romanian Liga I)
# ---------------------------------------------------------------------
$ws.Range("G4").Value = 6.8
$ws.Range("I4").Value = 1.75
$ws.Range("M4").Value = 1.06
$ws.Range("P4").Value = 1.95
$ws.Range("V4").Value = 2.32

# ---------------------------------------------------------------------
# 4) Insert a new row 5 for "Paraguayan Primera Division" fixture.
#    This shifts the former row 5 (Brazilian Serie A) down to row 6.
# ---------------------------------------------------------------------
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "Paraguayan Primera Division"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2025-10-16"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "18:30:00"
$ws.Range("D5").Value = "General Caballero"
$ws.Range("E5").Value = "Club Atletico Tembetary"
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 1.1
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 1.24
$ws.Range("Q5").Value = 1.01
$ws.Range("R5").Value = 1.18
$ws.Range("S5").Value = 1.02
$ws.Range("T5").Value = 1.03
$ws.Range("U5").Value = 1.03
$ws.Range("V5").Value = 1.01
$ws.Range("W5").Value = 1.01
$ws.Range("X5").Value = 1000
$ws.Range("Y5").Value = 1000
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 1000
$ws.Range("AD5").Value = 1000
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 1000
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

# ---------------------------------------------------------------------
# 5) Update row 6 (formerly row 5 - Brazilian Serie A) with new odds.
#    Text columns (League/Date/Time/Home/Away) already shifted down
#    correctly with the row insert, only numeric odds need updating.
# ---------------------------------------------------------------------
$ws.Range("F6").Value = 2.98
$ws.Range("G6").Value = 3.2
$ws.Range("H6").Value = 2.68
$ws.Range("I6").Value = 2.88
$ws.Range("J6").Value = 3.15
$ws.Range("K6").Value = 3.3
$ws.Range("L6").Value = 1.59
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 1.55
$ws.Range("O6").Value = 1.56
$ws.Range("P6").Value = 1.53
$ws.Range("Q6").Value = 2.6
$ws.Range("R6").Value = 1.18
$ws.Range("S6").Value = 5.6
$ws.Range("T6").Value = 1.94
$ws.Range("U6").Value = 1.67
$ws.Range("V6").Value = 1.53
$ws.Range("W6").Value = 1.45
$ws.Range("X6").Value = 8.800000000000001
$ws.Range("Y6").Value = 8.199999999999999
$ws.Range("Z6").Value = 980
$ws.Range("AA6").Value = 980
$ws.Range("AB6").Value = 8.800000000000001
$ws.Range("AC6").Value = 7.4
$ws.Range("AD6").Value = 13.5
$ws.Range("AE6").Value = 980
$ws.Range("AF6").Value = 22
$ws.Range("AG6").Value = 14.5
$ws.Range("AH6").Value = 980
$ws.Range("AI6").Value = 70
$ws.Range("AJ6").Value = 60
$ws.Range("AK6").Value = 980
$ws.Range("AL6").Value = 75
$ws.Range("AM6").Value = 230
$ws.Range("AN6").Value = 75
$ws.Range("AO6").Value = 980
